# "Generate Report for Handback"
# For each locale sheet (zh-cn, de-de):
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - A "Latest Target File" (col E) and "Latest Handback File" (col F) are now
#     populated (same file names/links as the handoff file + target xlf), and
#     "Latest Handback DateTime" (col G) gets a real timestamp.
# The Overview sheet mirrors the same Status text for each locale column.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet: just the status text changes ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus

function Update-LocaleSheet {
    param($SheetName, $Row2HandoffUrl, $Row2TargetUrl, $Row3HandoffUrl, $Row3TargetUrl, $HandbackDateTime)

    $ws = $wb.Worksheets.Item($SheetName)

    # Status -> handed back
    $ws.Range("B2").Value = $newStatus
    $ws.Range("B3").Value = $newStatus

    # Row 2 (c4d7e754... file): Latest Target File / Latest Handback File
    $fileNameE2 = $ws.Range("A2").Value2
    $fileNameF2 = $ws.Range("C2").Value2
    $ws.Range("E2").Value = $fileNameE2
    $ws.Range("F2").Value = $fileNameF2
    $ws.Hyperlinks.Add($ws.Range("E2"), $Row2HandoffUrl, "", "", $fileNameE2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $Row2TargetUrl, "", "", $fileNameF2) | Out-Null
    $ws.Range("G2").Value = $HandbackDateTime

    # Row 3 (ef03bb9d... file): Latest Target File / Latest Handback File
    $fileNameE3 = $ws.Range("A3").Value2
    $fileNameF3 = $ws.Range("C3").Value2
    $ws.Range("E3").Value = $fileNameE3
    $ws.Range("F3").Value = $fileNameF3
    $ws.Hyperlinks.Add($ws.Range("E3"), $Row3HandoffUrl, "", "", $fileNameE3) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $Row3TargetUrl, "", "", $fileNameF3) | Out-Null
    $ws.Range("G3").Value = $HandbackDateTime
}

Update-LocaleSheet `
    "zh-cn" `
    "https://github.com/OpenLocalizationTest/oltest/blob/d27144cb5f003b0291c8852cac18f26d84f9a402/e2e/c4d7e754-6531-44fa-aa41-a05b97856cc5.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d14c41b2cc2e7d0387d9bb97b4eb6121fcdd46d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/high/c4d7e754-6531-44fa-aa41-a05b97856cc5.e22def9ecd0f8241de21833c3bf31650fb85456d.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/d27144cb5f003b0291c8852cac18f26d84f9a402/e2e/ef03bb9d-4e78-4133-8858-a5cd822b575f.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d14c41b2cc2e7d0387d9bb97b4eb6121fcdd46d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/high/ef03bb9d-4e78-4133-8858-a5cd822b575f.01a0edaf24056f01c9d191c8a12634ce642f2338.zh-cn.xlf" `
    "2016-03-10 06:45:58"

Update-LocaleSheet `
    "de-de" `
    "https://github.com/OpenLocalizationTest/oltest/blob/d27144cb5f003b0291c8852cac18f26d84f9a402/e2e/c4d7e754-6531-44fa-aa41-a05b97856cc5.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3a52ac4c7f6a42504ebefd8d366573b4894089df/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/high/c4d7e754-6531-44fa-aa41-a05b97856cc5.e22def9ecd0f8241de21833c3bf31650fb85456d.de-de.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/d27144cb5f003b0291c8852cac18f26d84f9a402/e2e/ef03bb9d-4e78-4133-8858-a5cd822b575f.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3a52ac4c7f6a42504ebefd8d366573b4894089df/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/high/ef03bb9d-4e78-4133-8858-a5cd822b575f.01a0edaf24056f01c9d191c8a12634ce642f2338.de-de.xlf" `
    "2016-03-10 06:46:14"
